# Applies the edit described by the diff:
# 1. Slide 1: set showMasterSp="0" (DisplayMasterShapes = False)
# 2. Slide 1: remove the "Rectangle 2" shape (id=3), a white cover rectangle

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# 1. Hide/omit master shapes on this slide (adds showMasterSp="0" to <p:sld>)
$s.DisplayMasterShapes = 0

# 2. Delete the "Rectangle 2" shape
$s.Shapes.Item("Rectangle 2").Delete()
